$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared string "Start event" -> "Start" (cell B8)
$ws.Range("B8").Value = "Start"

# Update J8 and K8 values
$ws.Range("J8").Value = 2
$ws.Range("K8").Value = 10

# Update selection to L8
$ws.Range("L8").Select()
